# Update the marksheet's "Correct/Total" marks summary for roll 1401CS53.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Total correct answers (B11: Marking row's "Right" count placeholder -> corrected)
$ws.Range("B11").Value = 5

# Total marks scored (B12: "Total" row, "Right" column)
$ws.Range("B12").Value = 105

# Correct/total marks display text (E12: "Total" row, "Max" column)
$ws.Range("E12").Value = "105/140"
